# edit.ps1 - reproduce the authored changes to IndividualWorkSheet.xlsx
$wb = $excel.ActiveWorkbook
$nl = [char]10

# --- Sheet "신승민" (1st tab) : new task entries for row 2 and row 3 ---
$ws1 = $wb.Worksheets.Item(1)

# Write the new cell text in the same order the author originally typed it so the
# shared-string table is rebuilt with matching indices.
$ws1.Range("F2").Value = "디비 미완성으로 인해 알고리즘 수정 못함"
$ws1.Range("A2").Value = "web server구현 및 알고리즘 수정"
$ws1.Range("A3").Value = "web server구현 "
$ws1.Range("B3").Value = "부트스트랩 프레임워크를 이용한 웹페이지 UI 구축"
$ws1.Range("E2").Value = "1. lookupsmartphone.jsp(휴대폰 전제 보기 웹페이지 구현) -> 모빌라이저를 사용해 웹페이지 생성" + $nl + "2. lookupplan.jsp 구현 ->모빌라이저를 사용해 웹페이지 생성" + $nl + "3. recomsmartphon.jsp -> 휴대폰 추천 을 받았을때 띄워주는 창 구현   "
$ws1.Range("E3").Value = "1. detailphone.jsp 구현-> 휴대폰 상세 보기 기능( 이기능에 휴대폰의 스펙이나 성능을 설명해주는 동영상이 있으면 괜찮을것 같아서 youtube링크를 첨부하는쪽으로 추진) "

# Dates (stored as plain serials, not shared strings)
$ws1.Range("D2").Value = 43598
$ws1.Range("C3").Value = 43601
$ws1.Range("D3").Value = 43606

# Row heights grow to fit the new wrapped text
$ws1.Rows.Item(2).RowHeight = 134.4
$ws1.Rows.Item(3).RowHeight = 96

# Column widths widen for the new content
$ws1.Columns.Item(2).ColumnWidth = 16.857142857142858
$ws1.Columns.Item(5).ColumnWidth = 30.714285714285715
$ws1.Columns.Item(6).ColumnWidth = 23.857142857142858

# --- Sheet "박윤화" (4th tab): selection collapses to a single cell ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Activate()
$ws4.Range("A2").Select()

# --- Sheet "신승민" becomes the active tab, with the selection on E3 ---
$ws1.Activate()
$ws1.Range("E3").Select()
